# Update "想去人数" (want-to-go count) figures that changed in the
# upstream data refresh (gh-pages output regenerated at 456a3b4).
#
# Sheet 1 "展览"      (Worksheets.Item(1)) rows 2,3,12,14,15
# Sheet 4 "全部类型"  (Worksheets.Item(4)) rows 2,3,16,18,19
# are the same underlying events, so both sheets need the same updates.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item(1)   # 展览
$wsAll        = $wb.Worksheets.Item(4)   # 全部类型

# 展览 sheet
$wsExhibition.Range("F2").Value2  = 4538
$wsExhibition.Range("F3").Value2  = 2496
$wsExhibition.Range("F12").Value2 = 1698
$wsExhibition.Range("F14").Value2 = 3718
$wsExhibition.Range("F15").Value2 = 22

# 全部类型 sheet (mirrors the same events at different row numbers)
$wsAll.Range("F2").Value2  = 4538
$wsAll.Range("F3").Value2  = 2496
$wsAll.Range("F16").Value2 = 1698
$wsAll.Range("F18").Value2 = 3718
$wsAll.Range("F19").Value2 = 22
